$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.921.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.630.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('E3').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.522'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.44'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +10.04%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.865.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.633.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.571'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.69%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.92'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.40%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.947.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '9.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +19.05%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0707'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.29%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.32%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.58%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0494'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.92%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.30%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.430.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.31%  '
$ws.Range('E35').Style = 'Normal'

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.87'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.561'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.835'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.62%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0501'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.26%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '54.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '69.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.65%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.78%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.772.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '89.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.07%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0108'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.18%  '
$ws.Range('E51').Style = 'Normal'
